# Update cryptocurrency price/volume figures per the latest scrape run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.15"
$ws.Range("E2").Value = "'6.19%"
$ws.Range("D3").Value = "'32.48"
$ws.Range("E3").Value = "'10.70%"
$ws.Range("D4").Value = "'5.319"
$ws.Range("E4").Value = "'3.55%"
$ws.Range("D5").Value = "'0.07407"
$ws.Range("E5").Value = "'11.22%"
$ws.Range("D6").Value = "'7.742"
$ws.Range("E6").Value = "'5.43%"
$ws.Range("D7").Value = "'3.709"
$ws.Range("E7").Value = "'8.88%"
$ws.Range("D8").Value = "'1.594"
$ws.Range("E8").Value = "'17.86%"
$ws.Range("D9").Value = "'0.9231"
$ws.Range("E9").Value = "'0.52%"
$ws.Range("D10").Value = "'0.01632"
$ws.Range("E10").Value = "'2,420.15%"
$ws.Range("D11").Value = "'0.1669"
$ws.Range("E11").Value = "'6.30%"
$ws.Range("D12").Value = "'0.07381"
$ws.Range("E12").Value = "'12.59%"
$ws.Range("D13").Value = "'0.07982"
$ws.Range("E13").Value = "'4.18%"
$ws.Range("D14").Value = "'0.03113"
$ws.Range("E14").Value = "'7.07%"
$ws.Range("D15").Value = "'0.09820"
$ws.Range("E15").Value = "'9.22%"
$ws.Range("D16").Value = "'0.001532"
$ws.Range("E16").Value = "'-3.50%"
$ws.Range("D17").Value = "'0.04547"
$ws.Range("E17").Value = "'1.77%"
$ws.Range("D18").Value = "'0.006248"
$ws.Range("E18").Value = "'-0.11%"
$ws.Range("D19").Value = "'3.476"
$ws.Range("E19").Value = "'0.58%"
$ws.Range("D20").Value = "'2.238"
$ws.Range("E20").Value = "'0.45%"
$ws.Range("D21").Value = "'0.3271"
$ws.Range("E21").Value = "'1.81%"
$ws.Range("D22").Value = "'0.1311"
$ws.Range("E22").Value = "'0.16%"
$ws.Range("D23").Value = "'4.249"
$ws.Range("E23").Value = "'4.73%"
$ws.Range("D24").Value = "'0.1638"
$ws.Range("E24").Value = "'5.67%"
$ws.Range("E25").Value = "'2.81%"
$ws.Range("D26").Value = "'0.004532"
$ws.Range("E26").Value = "'9.79%"
$ws.Range("E27").Value = "'-6.35%"
$ws.Range("D28").Value = "'0.0001665"
$ws.Range("E28").Value = "'3.02%"
$ws.Range("D40").Value = "'0.04505"
$ws.Range("E40").Value = "'7.26%"
$ws.Range("D41").Value = "'0.007316"
$ws.Range("E41").Value = "'8.60%"
$ws.Range("D42").Value = "'0.1370"
$ws.Range("E42").Value = "'10.43%"
$ws.Range("D43").Value = "'0.002179"
$ws.Range("E43").Value = "'10.16%"
$ws.Range("D44").Value = "'0.01371"
$ws.Range("E44").Value = "'8.61%"
$ws.Range("D45").Value = "'0.00005958"
$ws.Range("E45").Value = "'4.05%"
$ws.Range("D46").Value = "'1.892"
$ws.Range("E46").Value = "'-3.83%"
$ws.Range("D47").Value = "'0.01299"
$ws.Range("E47").Value = "'-0.47%"
